$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) of the existing data rows down to the newly added rows 14-31,
# so the new "A" column numeric cells pick up the same style (border/bold/alignment) as rows 2-13.
for ($r = 14; $r -le 31; $r++) {
    $ws.Range("A13:F13").Copy()
    $ws.Range("A" + $r + ":F" + $r).PasteSpecial(-4122)
}
$excel.CutCopyMode = 0

# Write out the refreshed watchlist values (rows 2-31).
# Row 2
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "NSE:AMNPLST"
$ws.Range("C2").Value = "NSE:CASTROLIND"
$ws.Range("D2").Value = "NSE:ACC"
$ws.Range("F2").Value = "NSE:ABB"

# Row 3
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "NSE:APTECHT"
$ws.Range("C3").Value = "NSE:NAGREEKEXP"
$ws.Range("D3").Value = "NSE:ADANIENT"
$ws.Range("F3").Value = "NSE:ADANIPORTS"

# Row 4
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "NSE:BFUTILITIE"
$ws.Range("C4").Value = "NSE:PHOENIXLTD"
$ws.Range("D4").Value = "NSE:ADANIPORTS"
$ws.Range("F4").Value = "NSE:CUB"

# Row 5
$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "NSE:BIOFILCHEM"
$ws.Range("D5").Value = "NSE:AMBUJACEM"
$ws.Range("F5").Value = "NSE:IPCALAB"

# Row 6
$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "NSE:CAMS"
$ws.Range("D6").Value = "NSE:BAJFINANCE"
$ws.Range("F6").Value = "NSE:KOTAKBANK"

# Row 7
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "NSE:CONSUMBEES"
$ws.Range("D7").Value = "NSE:BOSCHLTD"
$ws.Range("F7").Value = "NSE:M&M"

# Row 8
$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "NSE:CPSEETF"
$ws.Range("D8").Value = "NSE:DALBHARAT"
$ws.Range("F8").Value = "NSE:MCDOWELL-N"

# Row 9
$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "NSE:CSLFINANCE"
$ws.Range("D9").Value = "NSE:DEEPAKNTR"
$ws.Range("F9").Value = "NSE:MGL"

# Row 10
$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "NSE:DCMSRIND"
$ws.Range("D10").Value = "NSE:INDIACEM"
$ws.Range("F10").Value = "NSE:OBEROIRLTY"

# Row 11
$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "NSE:DJML"
$ws.Range("D11").Value = "NSE:IOC"
$ws.Range("F11").Value = "NSE:PNB"

# Row 12
$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "NSE:EMUDHRA"
$ws.Range("D12").Value = "NSE:IPCALAB"
$ws.Range("F12").Value = "NSE:RELIANCE"

# Row 13
$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "NSE:FIEMIND"
$ws.Range("D13").Value = "NSE:JINDALSTEL"

# Row 14
$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "NSE:INFRABEES"
$ws.Range("D14").Value = "NSE:M&M"

# Row 15
$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "NSE:IPCALAB"
$ws.Range("D15").Value = "NSE:PIDILITIND"

# Row 16
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "NSE:ITDC"
$ws.Range("D16").Value = "NSE:POWERGRID"

# Row 17
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "NSE:JAIPURKURT"
$ws.Range("D17").Value = "NSE:RELIANCE"

# Row 18
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "NSE:JAYBARMARU"

# Row 19
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "NSE:JETFREIGHT"

# Row 20
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "NSE:JOCIL"

# Row 21
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "NSE:LAXMIMACH"

# Row 22
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "NSE:LT"

# Row 23
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "NSE:M&M"

# Row 24
$ws.Range("A24").Value = 22
$ws.Range("B24").Value = "NSE:MALLCOM"

# Row 25
$ws.Range("A25").Value = 23
$ws.Range("B25").Value = "NSE:MGL"

# Row 26
$ws.Range("A26").Value = 24
$ws.Range("B26").Value = "NSE:MOHEALTH"

# Row 27
$ws.Range("A27").Value = 25
$ws.Range("B27").Value = "NSE:MOHITIND"

# Row 28
$ws.Range("A28").Value = 26
$ws.Range("B28").Value = "NSE:MON100"

# Row 29
$ws.Range("A29").Value = 27
$ws.Range("B29").Value = "NSE:OBEROIRLTY"

# Row 30
$ws.Range("A30").Value = 28
$ws.Range("B30").Value = "NSE:RELIANCE"

# Row 31
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = "NSE:ROLLT"
